$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 289
$ws1.Range("G2").Value = 70
$ws1.Range("F4").Value = 1205
$ws1.Range("F5").Value = 837
$ws1.Range("F6").Value = 866
$ws1.Range("F7").Value = 1576
$ws1.Range("F9").Value = 1077
$ws1.Range("F10").Value = 37
$ws1.Range("F11").Value = 86
$ws1.Range("F12").Value = 211
$ws1.Range("F14").Value = 542
$ws1.Range("F15").Value = 84
$ws1.Range("F16").Value = 53
$ws1.Range("F20").Value = 597
$ws1.Range("F21").Value = 592
$ws1.Range("F22").Value = 76
$ws1.Range("F24").Value = 794
$ws1.Range("F25").Value = 270
$ws1.Range("F26").Value = 207

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1048
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F4").Value = 290
$ws2.Range("F6").Value = 206
$ws2.Range("F7").Value = 74
$ws2.Range("F9").Value = 96

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 273

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 273
$ws4.Range("F3").Value = 289
$ws4.Range("G3").Value = 70
$ws4.Range("F4").Value = 1048
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F6").Value = 1205
$ws4.Range("F7").Value = 837
$ws4.Range("F8").Value = 866
$ws4.Range("F9").Value = 1576
$ws4.Range("F11").Value = 1077
$ws4.Range("F12").Value = 37
$ws4.Range("F13").Value = 86
$ws4.Range("F14").Value = 211
$ws4.Range("F16").Value = 542
$ws4.Range("F17").Value = 84
$ws4.Range("F18").Value = 53
$ws4.Range("F21").Value = 290
$ws4.Range("F25").Value = 206
$ws4.Range("F26").Value = 206
$ws4.Range("F27").Value = 597
$ws4.Range("F28").Value = 592
$ws4.Range("F29").Value = 76
$ws4.Range("F31").Value = 794
$ws4.Range("F32").Value = 270
$ws4.Range("F33").Value = 74
$ws4.Range("F34").Value = 207
$ws4.Range("F36").Value = 96
$ws4.Range("F37").Value = 96
